$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 33
$ws.Range("F3").Value = 1
$ws.Range("F8").Value = 9473
$ws.Range("F9").Value = 238
$ws.Range("F10").Value = 33
$ws.Range("F11").Value = 689
$ws.Range("F12").Value = 1944
$ws.Range("F14").Value = 812
$ws.Range("F15").Value = 2599
$ws.Range("F17").Value = 3898
$ws.Range("F18").Value = 305
$ws.Range("F19").Value = 138
$ws.Range("F20").Value = 124
$ws.Range("F23").Value = 19
$ws.Range("F25").Value = 69
$ws.Range("F27").Value = 558
$ws.Range("F28").Value = 2
$ws.Range("F29").Value = 2137
$ws.Range("F30").Value = 1092
$ws.Range("F31").Value = 192
$ws.Range("F33").Value = 4315
$ws.Range("F35").Value = 163
$ws.Range("F36").Value = 331
$ws.Range("F37").Value = 144
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 187
$ws.Range("F3").Value = 973
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 187
$ws.Range("F3").Value = 33
$ws.Range("F4").Value = 973
$ws.Range("F5").Value = 1
$ws.Range("F11").Value = 9473
$ws.Range("F12").Value = 238
$ws.Range("F13").Value = 33
$ws.Range("F14").Value = 689
$ws.Range("F15").Value = 1944
$ws.Range("F17").Value = 812
$ws.Range("F19").Value = 2599
$ws.Range("F21").Value = 3898
$ws.Range("F22").Value = 305
$ws.Range("F24").Value = 124
$ws.Range("F27").Value = 19
$ws.Range("F30").Value = 69
$ws.Range("F32").Value = 558
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = 2137
$ws.Range("F35").Value = 1092
$ws.Range("F36").Value = 192
$ws.Range("F38").Value = 4315
$ws.Range("F40").Value = 163
$ws.Range("F41").Value = 331
$ws.Range("F42").Value = 144
